$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 142, shifting the existing rows 142-160 down to 143-161
$ws.Rows.Item(142).Insert()

# Populate the new row 142 with the weekly price entry
$ws.Range("A142").Value = 4
$ws.Range("B142").Value = 'Feria Lagunitas de Puerto Montt'
$ws.Range("C142").Value = 'Los Lagos'
$ws.Range("D142").Value = 44491
$ws.Range("E142").Value = 10
$ws.Range("F142").Value = 100112003
$ws.Range("G142").Value = 'Ajo'
$ws.Range("H142").Value = 'Chino'
$ws.Range("I142").Value = 'Primera'
$ws.Range("J142").Value = 280
$ws.Range("K142").Value = 20000
$ws.Range("L142").Value = 21000
$ws.Range("M142").Value = 20500
$ws.Range("N142").Value = '$/caja 10 kilos'
$ws.Range("O142").Value = 'China'
$ws.Range("P142").Value = 2050
$ws.Range("Q142").Value = 10
$ws.Range("R142").Value = 'Hortaliza'
